# Final Progress Pre Paper update (30-10-2023)
# Removes the two highest "Size" rows (14,15) for the Spinglass algorithm
# and refreshes the C (count) values for the remaining rows across all
# three algorithms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for Spinglass Size=14 and Size=15 (original rows 16-17).
# This shifts every following row up by two, which Excel also takes care
# of for the dimension reference automatically.
$ws.Rows("16:17").Delete()

# Refresh the remaining counts (column C) to their new values.
$ws.Range("C3").Value = 40
$ws.Range("C4").Value = 36
$ws.Range("C5").Value = 31
$ws.Range("C6").Value = 29
$ws.Range("C7").Value = 27
$ws.Range("C8").Value = 22
$ws.Range("C9").Value = 21
$ws.Range("C11").Value = 18
$ws.Range("C12").Value = 15
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 2
$ws.Range("C16").Value = 92
$ws.Range("C17").Value = 75
$ws.Range("C18").Value = 68
$ws.Range("C19").Value = 45
$ws.Range("C20").Value = 22
$ws.Range("C21").Value = 21
$ws.Range("C22").Value = 58
$ws.Range("C23").Value = 55
$ws.Range("C24").Value = 51
$ws.Range("C25").Value = 49
$ws.Range("C26").Value = 42
$ws.Range("C27").Value = 37
$ws.Range("C28").Value = 31
